$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.66"
$ws.Range("E2").Value = "'3.07%"
$ws.Range("D3").Value = "'39.52"
$ws.Range("E3").Value = "'3.33%"
$ws.Range("D4").Value = "'5.105"
$ws.Range("E4").Value = "'0.13%"
$ws.Range("D5").Value = "'0.08188"
$ws.Range("E5").Value = "'1.62%"
$ws.Range("D6").Value = "'2.027"
$ws.Range("E6").Value = "'4.75%"
$ws.Range("D7").Value = "'8.257"
$ws.Range("E7").Value = "'3.95%"
$ws.Range("D8").Value = "'0.9340"
$ws.Range("E8").Value = "'0.39%"
$ws.Range("D9").Value = "'0.1431"
$ws.Range("E9").Value = "'-0.16%"
$ws.Range("D10").Value = "'0.1984"
$ws.Range("E10").Value = "'3.45%"
$ws.Range("D11").Value = "'0.09137"
$ws.Range("E11").Value = "'1.80%"
$ws.Range("D12").Value = "'0.03553"
$ws.Range("E12").Value = "'1.18%"
$ws.Range("D13").Value = "'0.09824"
$ws.Range("E13").Value = "'0.52%"
$ws.Range("D14").Value = "'0.001403"
$ws.Range("E14").Value = "'0.31%"
$ws.Range("D15").Value = "'0.006354"
$ws.Range("E15").Value = "'3.40%"
$ws.Range("D16").Value = "'3.660"
$ws.Range("E16").Value = "'-1.79%"
$ws.Range("D17").Value = "'4.282"
$ws.Range("E17").Value = "'2.28%"
$ws.Range("D18").Value = "'3.291"
$ws.Range("E18").Value = "'-4.81%"
$ws.Range("D19").Value = "'0.3462"
$ws.Range("E19").Value = "'0.02%"
$ws.Range("E20").Value = "'-0.68%"
$ws.Range("D21").Value = "'4.827"
$ws.Range("E21").Value = "'0.15%"
$ws.Range("E22").Value = "'1.68%"
$ws.Range("D23").Value = "'0.04325"
$ws.Range("E23").Value = "'-0.59%"
$ws.Range("D24").Value = "'0.001224"
$ws.Range("E24").Value = "'-0.59%"
$ws.Range("D25").Value = "'0.004790"
$ws.Range("E25").Value = "'16.30%"
$ws.Range("E26").Value = "'-0.24%"
$ws.Range("D27").Value = "'0.0003998"
$ws.Range("E27").Value = "'-10.11%"
$ws.Range("D39").Value = "'0.02241"
$ws.Range("E39").Value = "'8.40%"
$ws.Range("D40").Value = "'0.05244"
$ws.Range("E40").Value = "'4.33%"
$ws.Range("D41").Value = "'0.007584"
$ws.Range("E41").Value = "'1.33%"
$ws.Range("D42").Value = "'0.009728"
$ws.Range("E42").Value = "'-3.89%"
$ws.Range("D43").Value = "'0.1380"
$ws.Range("E43").Value = "'2.45%"
$ws.Range("D44").Value = "'0.002139"
$ws.Range("E44").Value = "'-0.24%"
$ws.Range("D45").Value = "'0.009789"
$ws.Range("E45").Value = "'10.78%"
$ws.Range("D46").Value = "'0.00006400"
$ws.Range("E46").Value = "'3.19%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.23%"
$ws.Range("D48").Value = "'0.002767"
$ws.Range("E48").Value = "'-1.95%"
$ws.Range("E49").Value = "'-25.10%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.23%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.23%"
